$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Davis --- 12" (it also
# carries the _GoBack bookmark) and the paragraph right after it that
# starts "The deadlines for Honors students ...".
$davisPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Davis --- 12*") {
        $davisPara = $p
        break
    }
}

# Step 1: merge the "Davis --- 12" paragraph with the following
# paragraph by deleting its end-of-paragraph mark, while the "Davis"
# text is still present (deleting it now, before the text is removed,
# is what keeps the _GoBack bookmark alive through the merge).
$markRange = $d.Range($davisPara.Range.End - 1, $davisPara.Range.End)
$markRange.Delete()

# Step 2: remove the now-unwanted "Davis --- 12" text, leaving the
# bookmark and the rest of the (now combined) paragraph untouched.
$rng = $d.Content
$rng.Find.Execute("Davis --- 12", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$rng.Delete()

# Step 3: bold the paragraph mark of the combined paragraph (this is
# what shows up as <w:pPr><w:rPr><w:b/></w:rPr></w:pPr> in the XML).
# Bolding the whole paragraph range bolds both the mark and the visible
# run, so afterwards the run text is re-created clean (unbolded).
$mergedPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "The deadlines for Honors students*") {
        $mergedPara = $p
        break
    }
}
$mergedPara.Range.Font.Bold = $true

$textRng = $d.Content
$textRng.Find.Execute("The deadlines for Honors students to upload the " + `
    "final version of their thesis, including the signature cover " + `
    "sheet, to W&M Publish are:", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$insertStart = $textRng.Start
$savedText = $textRng.Text
$textRng.Delete()

$insPoint = $d.Range($insertStart, $insertStart)
$insPoint.Font.Bold = $false
$insPoint.InsertAfter($savedText)
